# Revert "Updated product Status"
# (this reverts commit 723b6c5f72226c25cc457578a32c78e3b4ccfd9e)
#
# The prior commit had filled in a "Status" value in column E (and B31)
# for a few task rows. Undo that by clearing those three cells back to
# empty, which also drops the now-unused shared strings ("done", "Jose",
# "Add booking form (order page) to front page?") on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E14").Value = $null
$ws.Range("E26").Value = $null
$ws.Range("B31").Value = $null

# Restore the view/selection state recorded in the reverted-to revision.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C19").Select()
